$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before DL; this shifts DL:MN (and all data in those
# columns, across every row) one column to the right, so old DL becomes DM,
# old DM becomes DN, ..., old MN becomes MO.
$ws.Columns("DL:DL").Insert()

# The insert above only carries over cells that already held a value/format;
# it leaves the brand-new column's cells on data rows unmaterialized. Copying
# from a neighboring genuinely-blank cell (DJ, which stays blank and
# untouched by the insert) stamps out a real-but-empty cell record at DL2/DL3,
# matching how the rest of the sheet already represents "empty" cells.
$ws.Range("DJ2:DJ3").Copy($ws.Range("DL2"))

# Set the header for the newly inserted column.
$ws.Range("DL1").Value = "DemonstrationProjectIdentifier"

# Update the Id values in the first data column (A2/A3).
$ws.Range("A2").Value = "6901488a7e79911955eafe38"
$ws.Range("A3").Value = "6901488a7e79911955eafe38"
